$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.743.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0628"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.871.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.664.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.525"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.763.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0507"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.292.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0175"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.819"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.795.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0525"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0972"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
